$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "menambahkan fitur kelas perkuliahan"
# The "Kode Prodi" / "Program Studi" columns (and their related helper text
# on the instructions sheet) are renamed to "Id_sms" / "Jenjang_pendidikan",
# and the instruction text is updated accordingly. The two header cells that
# back them on the "Petunjuk Pengisian" sheet get a red highlight fill so
# users don't edit them.
# ---------------------------------------------------------------------------

$idSms        = "Id_sms"
$jenjang      = "Jenjang_pendidikan"
$keteranganId = "Jangan diganti. Merupakan kode Prodi"
$keteranganJj = "Jangan diganti. Jejang Pendidikan Pada prodi tersebut"

# --- kelas ---
$wsKelas = $wb.Worksheets.Item("kelas")
$wsKelas.Range("F1").Value2 = $idSms
$wsKelas.Range("G1").Value2 = $jenjang

# --- krs ---
$wsKrs = $wb.Worksheets.Item("krs")
$wsKrs.Range("H1").Value2 = $idSms
$wsKrs.Range("I1").Value2 = $jenjang

# --- dosen ---
$wsDosen = $wb.Worksheets.Item("dosen")
$wsDosen.Range("J1").Value2 = $idSms
$wsDosen.Range("K1").Value2 = $jenjang

# --- nilai ---
$wsNilai = $wb.Worksheets.Item("nilai")
$wsNilai.Range("K1").Value2 = $idSms
$wsNilai.Range("L1").Value2 = $jenjang

# --- Petunjuk Pengisian ---
$wsHelp = $wb.Worksheets.Item("Petunjuk Pengisian")

$rows = @(8, 20, 33, 48)
foreach ($r in $rows) {
    $wsHelp.Range("B$r").Value2 = $idSms
    $wsHelp.Range("C$r").Value2 = $keteranganId
    $wsHelp.Range("B$r").Interior.Color = 255
}

$rows2 = @(9, 21, 34, 49)
foreach ($r in $rows2) {
    $wsHelp.Range("B$r").Value2 = $jenjang
    $wsHelp.Range("C$r").Value2 = $keteranganJj
    $wsHelp.Range("B$r").Interior.Color = 255
}

# ---------------------------------------------------------------------------
# View state: selections per sheet + which tab is active.
# ---------------------------------------------------------------------------

$wsKrs.Range("H1:I2").Select() | Out-Null

$wsDosen.Range("J1:K2").Select() | Out-Null

$wsNilai.Range("K1:L2").Select() | Out-Null

$wsHelp.Range("G43").Select() | Out-Null

$wsKelas.Activate()
$wsKelas.Range("G12").Select() | Out-Null
